$wb = $excel.ActiveWorkbook

# Update the calibrated grid battery capacity growth value (ReEDS calibration)
$dataSheet = $wb.Worksheets.Item("GBCGpUNR")
$dataSheet.Range("B2").Value = 2000

# Switch the active/selected sheet from "About" to "GBCGpUNR" and update its
# on-sheet selection to C10 (matches the saved view state in the edited file).
$dataSheet.Activate()
$dataSheet.Range("C10").Select()
